# Rename the sheet from Update_Labels to Add_Labels
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Add_Labels"

# Remove the now-unused D:E columns and the data rows (2:7), leaving only the header row
$ws.Range("D:E").EntireColumn.Delete()
$ws.Range("2:7").EntireRow.Delete()

# Update the header row content: Translation_id/Translation_value/Language/Status columns
# are gone; the remaining columns become Label_name, "1 US-en", "2 IN-kn"
$ws.Range("A1").Value = "Label_name"
$ws.Range("B1").Value = "1 US-en"
$ws.Range("C1").Value = "2 IN-kn"

# Column widths: A=28, B:C=20 (COM ColumnWidth excludes the ~5px/6 gridline
# padding that OOXML's stored `width` includes, so add 5/6 before assigning)
$ws.Columns.Item(1).ColumnWidth = 27.166666666666668
$ws.Columns.Item(2).ColumnWidth = 19.166666666666668
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
